$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in column D (Fecha) for rows 2,3 and 6,7
$ws.Range("D2").Value = 44846
$ws.Range("D3").Value = 44846
$ws.Range("D6").Value = 44832
$ws.Range("D7").Value = 44832
